# Correcting Relevance Markers Appenzeller-Herzog (2019) - van Dis (2020)
# Updates record_atd / td_sim_1 (column C) values, and the derived average (column D)
# for the rows whose simulated time-to-discovery changed, plus the recomputed
# average_simulation_TD summary cell (C170).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 3; C = 3; D = 3 }
    @{ Row = 5; C = 412; D = 412 }
    @{ Row = 7; C = 640; D = 640 }
    @{ Row = 9; C = 345; D = 345 }
    @{ Row = 11; C = 66; D = 66 }
    @{ Row = 13; C = 119; D = 119 }
    @{ Row = 15; C = 51; D = 51 }
    @{ Row = 17; C = 208; D = 208 }
    @{ Row = 19; C = 148; D = 148 }
    @{ Row = 21; C = 225; D = 225 }
    @{ Row = 23; C = 14; D = 14 }
    @{ Row = 25; C = 28; D = 28 }
    @{ Row = 27; C = 80; D = 80 }
    @{ Row = 29; C = 15; D = 15 }
    @{ Row = 32; C = 121; D = 124 }
    @{ Row = 34; C = 120; D = 120 }
    @{ Row = 35; C = 86; D = 86 }
    @{ Row = 37; C = 316; D = 316 }
    @{ Row = 39; C = 214; D = 214 }
    @{ Row = 41; C = 16; D = 16 }
    @{ Row = 43; C = 83; D = 108.5 }
    @{ Row = 44; C = 125; D = 125 }
    @{ Row = 46; C = 44; D = 44 }
    @{ Row = 48; C = 81; D = 81 }
    @{ Row = 50; C = 23; D = 23 }
    @{ Row = 52; C = 95; D = 95 }
    @{ Row = 55; C = 18; D = 18 }
    @{ Row = 56; C = 101; D = 101 }
    @{ Row = 58; C = 162; D = 162 }
    @{ Row = 60; C = 209; D = 209 }
    @{ Row = 62; C = 88; D = 88 }
    @{ Row = 67; C = 118; D = 118 }
    @{ Row = 68; C = 270; D = 270 }
    @{ Row = 70; C = 141; D = 141 }
    @{ Row = 72; C = 133; D = 73.5 }
    @{ Row = 74; C = 46; D = 46 }
    @{ Row = 75; C = 231; D = 231 }
    @{ Row = 77; C = 93; D = 93 }
    @{ Row = 79; C = 236; D = 236 }
    @{ Row = 81; C = 56; D = 56 }
    @{ Row = 83; C = 388; D = 388 }
    @{ Row = 85; C = 6; D = 12.5 }
    @{ Row = 87; C = 11; D = 11 }
    @{ Row = 88; C = 227; D = 227 }
    @{ Row = 90; C = 216; D = 216 }
    @{ Row = 91; C = 206; D = 206 }
    @{ Row = 93; C = 30; D = 30 }
    @{ Row = 95; C = 201; D = 201 }
    @{ Row = 97; C = 219; D = 219 }
    @{ Row = 99; C = 36; D = 36 }
    @{ Row = 101; C = 167; D = 167 }
    @{ Row = 103; C = 349; D = 349 }
    @{ Row = 105; C = 7; D = 7 }
    @{ Row = 107; C = 274; D = 274 }
    @{ Row = 110; C = 458; D = 458 }
    @{ Row = 111; C = 52; D = 122.5 }
    @{ Row = 113; C = 215; D = 215 }
    @{ Row = 114; C = 109; D = 239.5 }
    @{ Row = 116; C = 365; D = 365 }
    @{ Row = 117; C = 94; D = 74.5 }
    @{ Row = 119; C = 35; D = 35 }
    @{ Row = 120; C = 271; D = 271 }
    @{ Row = 122; C = 10; D = 82 }
    @{ Row = 123; C = 161; D = 161 }
    @{ Row = 125; C = 73; D = 73 }
    @{ Row = 127; C = 288; D = 288 }
    @{ Row = 129; C = 50; D = 50 }
    @{ Row = 131; C = 276; D = 276 }
    @{ Row = 133; C = 41; D = 41 }
    @{ Row = 135; C = 13; D = 13 }
    @{ Row = 137; C = 29; D = 29 }
    @{ Row = 139; C = 38; D = 38 }
    @{ Row = 141; C = 427; D = 427 }
    @{ Row = 143; C = 107; D = 107 }
    @{ Row = 145; C = 164; D = 164 }
    @{ Row = 147; C = 246; D = 246 }
    @{ Row = 149; C = 302; D = 302 }
    @{ Row = 152; C = 484; D = 484 }
    @{ Row = 153; C = 123; D = 123 }
    @{ Row = 155; C = 77; D = 77 }
    @{ Row = 157; C = 341; D = 341 }
    @{ Row = 159; C = 33; D = 33 }
    @{ Row = 161; C = 12; D = 12 }
    @{ Row = 163; C = 308; D = 308 }
    @{ Row = 165; C = 399; D = 399 }
    @{ Row = 167; C = 76; D = 76 }
    @{ Row = 169; C = 89; D = 89 }
    @{ Row = 170; C = 155.4157303370787 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    if ($u.ContainsKey("D")) {
        $ws.Cells.Item($u.Row, 4).Value = $u.D
    }
}

